$wb = $excel.ActiveWorkbook

# --- 1. Refresh the query timestamps on the existing "data" sheet (col F) ---
$dataSheet = $wb.Worksheets.Item("data")
$dataSheet.Range("F2").Value = "2021-10-05 14:22:27.845125"
$dataSheet.Range("F3").Value = "2021-10-05 14:22:27.845130"
$dataSheet.Range("F4").Value = "2021-10-05 14:22:27.845132"
$dataSheet.Range("F5").Value = "2021-10-05 14:22:27.845134"
$dataSheet.Range("F6").Value = "2021-10-05 14:22:27.845136"

# --- 2. Add a new "metadata" sheet right after "data" ---
$metaSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dataSheet)
$metaSheet.Name = "metadata"

# Header row values
$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

# Apply the same bold/centered/bordered header style used on "data"!B1:F1
$dataSheet.Range("B1:F1").Copy()
$metaSheet.Range("B1:G1").PasteSpecial(-4122)  # xlPasteFormats

# Data row values
$metaSheet.Range("A2").Value = 0
$metaSheet.Range("B2").Value = "Refuted genes"
$metaSheet.Range("C2").Value = 8
$metaSheet.Range("E2").Value = "2019-07-09T14:10:48.357036Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:22:27.842775"
$metaSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/8/?format=json"

# data_version ("0.13") must stay textual, not become the number 0.13
$metaSheet.Range("D2").NumberFormat = "@"
$metaSheet.Range("D2").Value = "0.13"
$metaSheet.Range("D2").Style = "Normal"

# A2 carries the same style as "data"!A2 (bold/centered/bordered, numeric 0)
$dataSheet.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)  # xlPasteFormats

$dataSheet.Activate()
